$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5433
$ws1.Range("F4").Value = 33
$ws1.Range("F7").Value = 627
$ws1.Range("F8").Value = 599
$ws1.Range("F9").Value = 1062
$ws1.Range("F11").Value = 1500
$ws1.Range("F12").Value = 4648
$ws1.Range("F13").Value = 445
$ws1.Range("F14").Value = 203
$ws1.Range("F15").Value = 179
$ws1.Range("F16").Value = 100
$ws1.Range("F17").Value = 3569
$ws1.Range("F18").Value = 184
$ws1.Range("F19").Value = 1122
$ws1.Range("F20").Value = 108
$ws1.Range("F23").Value = 28
$ws1.Range("F24").Value = 140
$ws1.Range("F27").Value = 75
$ws1.Range("F28").Value = 324
$ws1.Range("F29").Value = 35
$ws1.Range("F31").Value = 21
$ws1.Range("F32").Value = 31
$ws1.Range("F33").Value = 35
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5433
$ws4.Range("F5").Value = 33
$ws4.Range("F8").Value = 627
$ws4.Range("F9").Value = 599
$ws4.Range("F10").Value = 1062
$ws4.Range("F12").Value = 1500
$ws4.Range("F13").Value = 4648
$ws4.Range("F14").Value = 445
$ws4.Range("F15").Value = 203
$ws4.Range("F16").Value = 179
$ws4.Range("F17").Value = 100
$ws4.Range("F18").Value = 3569
$ws4.Range("F19").Value = 184
$ws4.Range("F20").Value = 1122
$ws4.Range("F21").Value = 108
$ws4.Range("F24").Value = 28
$ws4.Range("F25").Value = 140
$ws4.Range("F28").Value = 75
$ws4.Range("F29").Value = 324
$ws4.Range("F30").Value = 35
$ws4.Range("F32").Value = 21
$ws4.Range("F33").Value = 31
$ws4.Range("F34").Value = 35
